$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: STATUS_SISWA ---
$ws.Range("J1").Value = "STATUS_SISWA"
$ws.Range("J2").Value = "TIDAK AKTIF"
$ws.Range("J3").Value = "AKTIF"
$ws.Range("J4").Value = "AKTIF"
$ws.Range("J5").Value = "TIDAK AKTIF"
$ws.Range("J6").Value = "AKTIF"
$ws.Range("J7").Value = "AKTIF"

# --- Shift the L-column instruction notes down to make room for new STATUS_SISWA note ---
$ws.Range("L8").Insert(-4121)
$ws.Range("L7").Copy($ws.Range("L8"))
$ws.Range("L7").ClearContents()
$ws.Range("L7").Value = "*untuk format isian kolom STATUS_SISWA wajib diisi (TIDAK AKTIF/AKTIF) (WAJIB DI ISI)"
$ws.Range("L7").Characters(26, 14).Font.Bold = $true
$ws.Range("L7").Characters(26, 14).Font.Color = 255
$ws.Range("L7").Characters(40, 31).Font.Color = 255
$ws.Range("L7").Characters(71, 15).Font.Bold = $true
$ws.Range("L7").Characters(71, 15).Font.Color = 255

# --- Header row (row 1): center/center alignment for all; A,B,C,G,H,I,J keep bold-red font, D,F,E keep bold font ---
$ws.Range("A1:J1").HorizontalAlignment = -4108
$ws.Range("A1:J1").VerticalAlignment = -4108

# --- Body cells A2:J7: center/center alignment, default font ---
$ws.Range("A2:D7").HorizontalAlignment = -4108
$ws.Range("A2:D7").VerticalAlignment = -4108
$ws.Range("F2:J7").HorizontalAlignment = -4108
$ws.Range("F2:J7").VerticalAlignment = -4108

# --- Column E (NAMA_WALI) body cells: center/center, text numfmt already applied ---
$ws.Range("E2:E7").HorizontalAlignment = -4108
$ws.Range("E2:E7").VerticalAlignment = -4108

# --- L column notes: left/center alignment ---
$ws.Range("L3:L9").HorizontalAlignment = -4131
$ws.Range("L3:L9").VerticalAlignment = -4108
